# Estado de Cuenta (NIT-9009938634): refresh the "Periodo Mora" column.
# Previous statement periods (2012, 2101, 2102, 2103) are replaced with the
# updated set, written in reverse chronological order (2103, 2102, 2101, 2012).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16").Value = "2103"
$ws.Range("E17").Value = "2102"
$ws.Range("E18").Value = "2101"
$ws.Range("E19").Value = "2012"
